$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old A7 entry (List_of_Atlas_launches_(2010-2019) link), which
# drops the shared string out of the sheet so the sst gets renumbered.
$ws.Range("A7").ClearContents()

# Add new column B with "Thor and Delta" launch list links.
$ws.Range("B1").Value = "https://en.wikipedia.org/wiki/List_of_Thor_and_Delta_launches_(1957%E2%80%9359)"
$ws.Range("B2").Value = "https://en.wikipedia.org/wiki/List_of_Thor_and_Delta_launches_(1960%E2%80%9369)"
$ws.Range("B3").Value = "https://en.wikipedia.org/wiki/List_of_Thor_and_Delta_launches_(1970%E2%80%9379)"
$ws.Range("B4").Value = "https://en.wikipedia.org/wiki/List_of_Thor_and_Delta_launches_(1980%E2%80%9389)"
$ws.Range("B5").Value = "https://en.wikipedia.org/wiki/List_of_Thor_and_Delta_launches_(1990%E2%80%9399)"
$ws.Range("B6").Value = "https://en.wikipedia.org/wiki/List_of_Thor_and_Delta_launches_(2000%E2%80%9309)"

# Match the new column B width (72 characters, bestFit-style) used for the URLs.
$ws.Columns.Item(2).ColumnWidth = 71.17

# Selection moves to the now-empty A7, matching the saved view state.
$ws.Range("A7").Select()
